$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (rows 24-27)
$ws.Range("A24").Value = 23
$ws.Range("C24").Value = 0.07
$ws.Range("D24").Value = 0.1
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = "None"
$ws.Range("G24").Value = "No"

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "US"
$ws.Range("C25").Value = 0.07
$ws.Range("D25").Value = 0.1
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = "None"
$ws.Range("G25").Value = "No"

$ws.Range("A26").Value = 25
$ws.Range("C26").Value = 0.07
$ws.Range("D26").Value = 0.1
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = "Yes"
$ws.Range("G26").Value = "No"

$ws.Range("A27").Value = 26
$ws.Range("C27").Value = 0.0626
$ws.Range("D27").Value = 0.1
$ws.Range("E27").Value = 40
$ws.Range("F27").Value = "Yes"
$ws.Range("G27").Value = "No"
$ws.Range("I27").Value = "Add TE to MV"

# Column G width (closest achievable value to the target 11.140625 given this
# runtime's pixel-snapped ColumnWidth rounding)
$ws.Columns.Item(7).ColumnWidth = 10.333333333333332

# Sheet view changes: scroll so row 4 is at the top, then select H27
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("H27").Select()

# Window position/size
$win.Left = 31200
$win.Top = 0
$win.Width = 17100
$win.Height = 17400
